$d = $word.ActiveDocument

# Namespace-qualified OOXML fragment wrapper used with Range.InsertXML.
# When InsertXML is invoked on a RANGE THAT IS COLLAPSED AT THE START of an
# existing paragraph, the block content in the fragment is merged in as a
# prefix of that paragraph (its own trailing paragraph mark is discarded and
# fused with the target paragraph) - this lets us inject multiple distinct
# <w:r>/<w:proofErr> children into a paragraph without Word's usual "merge
# same-formatted runs" behavior collapsing them back into one <w:t>.
function New-Pkg([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1) Drop the "Hinbekommen, das man sich in eine Fahrgemeinschaft eintragen
#    kann" paragraph entirely - superseded by the password-change note below.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Hinbekommen*") { $target = $i; break }
}
if ($target -ne $null) {
    $d.Paragraphs.Item($target).Range.Delete()
}

# 2) "Mit derselben Methode ... Passwort ändern" gains a second run
#    "(Name eintragen)" appended in its own <w:r> (not merged into the
#    first). Replace the paragraph's text with an XML fragment holding both
#    runs, inserted at the paragraph's start so it fuses in cleanly and
#    keeps the paragraph's own identity/properties.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Mit derselben Methode*") { $target = $i; break }
}
$p = $d.Paragraphs.Item($target)
$rng = $p.Range
$textLen = ($rng.End - $rng.Start) - 1
$textOnly = $d.Range($rng.Start, $rng.Start + $textLen)
$textOnly.Delete()
$p = $d.Paragraphs.Item($target)
$insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
$fragXml = "<w:p><w:r><w:t>Mit derselben Methode in Einstellungen Name und Passwort ändern</w:t></w:r><w:r><w:t>(Name eintragen)</w:t></w:r></w:p>"
$insertPoint.InsertXML((New-Pkg $fragXml))

# 3) Insert a brand-new paragraph "Erinnnerung wenn das Auto voll ist und
#    wenn kein Platz mehr ist" right before "Fehlerbehbung". "Erinnnerung"
#    is flagged as a (mis-)spelling via proofErr start/end markers, as Word
#    would do while the author was typing it. Insert at the start of the
#    "Fehlerbehbung" paragraph plus an extra trailing empty <w:p/> so the
#    fragment lands as its own, separate paragraph ahead of it, then drop
#    the now-spurious empty paragraph that InsertXML leaves behind.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Fehlerbehbung*") { $target = $i; break }
}
$p = $d.Paragraphs.Item($target)
$insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
$fragXml = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Erinnnerung</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> wenn das Auto voll ist und wenn kein Platz mehr ist</w:t></w:r></w:p><w:p/>'
$insertPoint.InsertXML((New-Pkg $fragXml))
# the paragraph right after our new one is the spurious empty one
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Erinnnerung*") { $target2 = $i + 1; break }
}
if ($target2 -ne $null -and $d.Paragraphs.Item($target2).Range.Text.Length -le 1) {
    $d.Paragraphs.Item($target2).Range.Delete()
}

# 4) Flag "Fehlerbehbung" itself as a (mis-)spelling too, wrapping the
#    existing run in proofErr start/end markers.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Fehlerbehbung*") { $target = $i; break }
}
$p = $d.Paragraphs.Item($target)
$rng = $p.Range
$textLen = ($rng.End - $rng.Start) - 1
$textOnly = $d.Range($rng.Start, $rng.Start + $textLen)
$textOnly.Delete()
$p = $d.Paragraphs.Item($target)
$insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
$fragXml = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Fehlerbehbung</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$insertPoint.InsertXML((New-Pkg $fragXml))

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "Para ${i}: [$($d.Paragraphs.Item($i).Range.Text)]"
}
